$wb = $excel.ActiveWorkbook

# Update both sheets that carry the "展览" style listing data:
# - "展览" (exhibitions)
# - "全部类型" (all types) which mirrors the same rows

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 1465
    $ws.Range("F9").Value = 256
}
